$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: extend the header sequence with new P1=14, Q1=15 cells,
# copying the bold/centered/bordered format from the existing O1 header cell.
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# Rows 2-25: swap the I/K values and the M/O values, then append
# two new data columns P and Q (both valued 2) to each row.
for ($r = 2; $r -le 25; $r++) {
    $ws.Range("I$r").Value = 2
    $ws.Range("K$r").Value = 1
    $ws.Range("M$r").Value = 2
    $ws.Range("O$r").Value = 1
    $ws.Range("P$r").Value = 2
    $ws.Range("Q$r").Value = 2
}

Write-Output "edit applied"
